# Updated cryptos list on Fri Jun  7 19:55:27 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.365.73"
$ws.Range("E2").Value = "  -1.75%  "

$ws.Range("D3").Value = "'3.692.89"
$ws.Range("E3").Value = "  -2.72%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'682.25"
$ws.Range("E5").Value = "  -3.26%  "

$ws.Range("D6").Value = "'162.30"
$ws.Range("E6").Value = "  -3.87%  "

$ws.Range("D7").Value = "'3.689.94"
$ws.Range("E7").Value = "  -2.65%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  -5.02%  "

$ws.Range("D10").Value = "'0.148"
$ws.Range("E10").Value = "  -7.04%  "

$ws.Range("D11").Value = "'7.26"
$ws.Range("E11").Value = "  -0.93%  "

$ws.Range("E12").Value = "  -0.75%  "

$ws.Range("D13").Value = "'0.0000236"
$ws.Range("E13").Value = "  -6.26%  "

$ws.Range("D14").Value = "'33.53"
$ws.Range("E14").Value = "  -7.09%  "

$ws.Range("D15").Value = "'4.317.83"
$ws.Range("E15").Value = "  -2.73%  "

$ws.Range("D16").Value = "'3.684.62"
$ws.Range("E16").Value = "  -4.30%  "

$ws.Range("D17").Value = "'69.414.54"
$ws.Range("E17").Value = "  -1.72%  "

$ws.Range("E18").Value = "  -1.53%  "

$ws.Range("D19").Value = "'16.33"
$ws.Range("E19").Value = "  -5.11%  "

$ws.Range("D20").Value = "'6.60"
$ws.Range("E20").Value = "  -7.44%  "

$ws.Range("D21").Value = "'485.26"
$ws.Range("E21").Value = "  -0.98%  "

$ws.Range("D22").Value = "'9.80"
$ws.Range("E22").Value = "  -6.64%  "

$ws.Range("D23").Value = "'0.663"
$ws.Range("E23").Value = "  -8.38%  "

$ws.Range("D24").Value = "'79.68"
$ws.Range("E24").Value = "  -6.26%  "

$ws.Range("D25").Value = "'3.839.48"
$ws.Range("E25").Value = "  -2.83%  "

$ws.Range("D26").Value = "'0.0000128"
$ws.Range("E26").Value = "  -10.42%  "

$ws.Range("D27").Value = "'11.59"
$ws.Range("E27").Value = "  -3.61%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").Value = "'9.54"
$ws.Range("E29").Value = "  -8.33%  "

$ws.Range("D30").Value = "'1.81"
$ws.Range("E30").Value = "  -11.79%  "

$ws.Range("E31").Value = "  -10.62%  "

$ws.Range("E32").Value = "  -4.21%  "

$ws.Range("D33").Value = "'6.72"
$ws.Range("E33").Value = "  -8.07%  "

$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("D35").Value = "'26.80"
$ws.Range("E35").Value = "  -7.67%  "

$ws.Range("E36").Value = "  -4.78%  "

$ws.Range("D37").Value = "'3.662.40"
$ws.Range("E37").Value = "  -2.78%  "

$ws.Range("D38").Value = "'8.51"
$ws.Range("E38").Value = "  -5.54%  "

$ws.Range("D39").Value = "'6.07"
$ws.Range("E39").Value = "  +3.08%  "

$ws.Range("E40").Value = "  -7.02%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("E42").Value = "  -4.44%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").Value = "'0.958"
$ws.Range("E44").Value = "  -7.76%  "

$ws.Range("D45").Value = "'159.79"
$ws.Range("E45").Value = "  -2.74%  "

$ws.Range("D46").Value = "'48.24"
$ws.Range("E46").Value = "  -0.95%  "

$ws.Range("D47").Value = "'2.82"
$ws.Range("E47").Value = "  -13.12%  "

$ws.Range("D48").Value = "'393.11"
$ws.Range("E48").Value = "  -6.57%  "

$ws.Range("D49").Value = "'0.000277"
$ws.Range("E49").Value = "  -9.38%  "

$ws.Range("D50").Value = "'1.29"
$ws.Range("E50").Value = "  -3.97%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'8.04"
$ws.Range("E51").Value = "  -7.24%  "
